# Add a new worksheet "bead_size_t_test" with bead size measurement data
$wb = $excel.ActiveWorkbook

# Add a new sheet at the end of the workbook (after the current last sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "bead_size_t_test"

# Header row
$newSheet.Range("A1").Value = "measurement"
$newSheet.Range("B1").Value = "size"

# Data rows
$data = @(
    @(6.24, "large"),
    @(5.4, "large"),
    @(4.78, "large"),
    @(3.89, "small"),
    @(4.21, "small"),
    @(4.13, "small"),
    @(3.68, "small"),
    @(3.4, "small")
)

$row = 2
foreach ($item in $data) {
    $newSheet.Cells.Item($row, 1).Value = $item[0]
    $newSheet.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# Auto-fit column A to its content (matches the "bestFit" width behavior)
$newSheet.Columns.Item(1).AutoFit() | Out-Null

# Select full column A on the new sheet, and make it the active/selected sheet
$newSheet.Activate()
$newSheet.Range("A1:A1048576").Select()
